# Adding test cases for watch list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 33 (E33): SKIP -> PASS
$ws.Range("E33").Value = "PASS"

# Copy formatting of the last existing data row (49) down into the two
# new rows so the new cells inherit the same borders/fonts/styles.
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E51").PasteSpecial(-4122)

# New row 50 - FollowUnfollowPostsAuthor / OPQA-427
$ws.Range("A50").Value = "FollowUnfollowPostsAuthor"
$ws.Range("B50").Value = "OPQA-427"
$ws.Range("C50").Value = "Veirfy that the user is able to follow the author of the post directly from the post"
$ws.Range("D50").Value = "Y"
$ws.Range("E50").Value = "PASS"

# New row 51 - CommentOnUsersOwnPost / OPQA-377
$ws.Range("A51").Value = "CommentOnUsersOwnPost"
$ws.Range("B51").Value = "OPQA-377"
$ws.Range("C51").Value = "Verify that the user is able to comment on the post a user authored themselves."
$ws.Range("D51").Value = "Y"
$ws.Range("E51").Value = "PASS"

# Add a hyperlink for the new Jira reference in B50, matching the pattern
# used by all the other JIRA-ID cells in column B.
$opqa427 = "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-427"
$ws.Hyperlinks.Add($ws.Range("B50"), $opqa427, "", "", $opqa427)
$ws.Range("B50").Value = "OPQA-427"
$ws.Range("B50").Font.Underline = 2
$ws.Range("B50").Font.Underline = 2

# Keep the active selection on the new last row, like Excel would after
# appending rows at the bottom of the table.
$ws.Range("A51:E51").Select()
